$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as TEXT, preserving the original (unstyled) cell
# format so numeric-looking strings (e.g. "591.57") are not silently
# auto-converted into floating point numbers / dates by Excel's input
# parser. We flip the cell to the "@" (Text) number format just long enough
# to assign the literal string, then restore the cell's original Style so
# no visible formatting changes leak into the saved file.
function Set-TextValue {
    param($cellRef, $val)
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

$ws.Range('D2').Value = '64.555.45'
$ws.Range('E2').Value = '  +1.39%  '
$ws.Range('D3').Value = '3.157.57'
$ws.Range('E3').Value = '  +0.91%  '
$ws.Range('E4').Value = '  -0.05%  '
Set-TextValue 'D5' '591.57'
$ws.Range('E5').Value = '  +0.18%  '
Set-TextValue 'D6' '147.31'
$ws.Range('E6').Value = '  +0.83%  '
$ws.Range('D8').Value = '3.154.90'
$ws.Range('E8').Value = '  +1.09%  '
$ws.Range('E9').Value = '  -0.74%  '
Set-TextValue 'D10' '0.162'
$ws.Range('E10').Value = '  -0.27%  '
Set-TextValue 'D11' '5.97'
$ws.Range('E11').Value = '  +4.74%  '
Set-TextValue 'D12' '0.463'
$ws.Range('E12').Value = '  -0.99%  '
$ws.Range('E13').Value = '  -1.94%  '
Set-TextValue 'D14' '37.23'
$ws.Range('E14').Value = '  +3.01%  '
$ws.Range('D15').Value = '3.679.74'
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('E16').Value = '  -1.14%  '
Set-TextValue 'D17' '7.22'
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').Value = '64.192.18'
$ws.Range('E18').Value = '  +0.94%  '
$ws.Range('D19').Value = '3.152.49'
$ws.Range('E19').Value = '  +0.85%  '
Set-TextValue 'D20' '469.11'
$ws.Range('E20').Value = '  +0.61%  '
Set-TextValue 'D21' '14.48'
$ws.Range('E21').Value = '  +1.67%  '
Set-TextValue 'D22' '0.737'
$ws.Range('E22').Value = '  +0.11%  '
Set-TextValue 'D23' '7.52'
$ws.Range('E23').Value = '  -0.29%  '
$ws.Range('B24').Value = 'Fetch.AI'
$ws.Range('C24').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D24' '2.35'
$ws.Range('E24').Value = '  +8.39%  '
$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D25' '13.08'
$ws.Range('E25').Value = '  -1.52%  '
Set-TextValue 'D26' '81.55'
$ws.Range('E26').Value = '  -0.85%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('E28').Value = '  +11.52%  '
Set-TextValue 'D29' '2.73'
$ws.Range('E29').Value = '  +0.64%  '
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D30' '7.38'
$ws.Range('E30').Value = '  +8.11%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D31' '2.23'
$ws.Range('E31').Value = '  +0.48%  '
$ws.Range('E32').Value = '  -0.14%  '
Set-TextValue 'D33' '27.57'
$ws.Range('E33').Value = '  +1.72%  '
$ws.Range('E34').Value = '  +0.61%  '
$ws.Range('D35').Value = '0.0₃0852'
$ws.Range('E35').Value = '  -1.08%  '
$ws.Range('E36').Value = '  +1.28%  '
Set-TextValue 'D37' '2.35'
$ws.Range('E37').Value = '  -1.73%  '
$ws.Range('E38').Value = '  -0.15%  '
Set-TextValue 'D39' '3.29'
$ws.Range('E39').Value = '  -2.06%  '
Set-TextValue 'D40' '51.98'
$ws.Range('E40').Value = '  +2.16%  '
Set-TextValue 'D41' '457.38'
$ws.Range('E41').Value = '  +2.21%  '
$ws.Range('E42').Value = '  +4.33%  '
Set-TextValue 'D43' '0.295'
$ws.Range('E43').Value = '  +6.12%  '
Set-TextValue 'D44' '0.0374'
$ws.Range('E44').Value = '  +0.56%  '
$ws.Range('D45').Value = '2.935.59'
$ws.Range('E45').Value = '  +0.50%  '
Set-TextValue 'D46' '40.56'
$ws.Range('E46').Value = '  +15.32%  '
Set-TextValue 'D47' '0.110'
$ws.Range('E47').Value = '  -0.81%  '
Set-TextValue 'D48' '128.20'
$ws.Range('E48').Value = '  +2.01%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('E50').Value = '  +2.63%  '
$ws.Range('E51').Value = '  -0.23%  '
